$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("I3").Value = 6
$ws.Range("Y3").Value = 1.53
$ws.Range("Z3").Value = 2.38
$ws.Range("AO3").Value = 26
$ws.Range("AP3").Value = 19
$ws.Range("AS3").Value = 51

# Row 4
$ws.Range("O4").Value = 1.44
$ws.Range("P4").Value = 2.75
$ws.Range("AF4").Value = 9

# Row 5
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 7

# Row 6
$ws.Range("I6").Value = 2.3
$ws.Range("AQ6").Value = 21

# Row 8
$ws.Range("G8").Value = 3.8
$ws.Range("J8").Value = 4.2
$ws.Range("K8").Value = 2.15
$ws.Range("N8").Value = 6.7
$ws.Range("O8").Value = 1.35
$ws.Range("P8").Value = 2.92
$ws.Range("S8").Value = 2.05
$ws.Range("T8").Value = 1.7
$ws.Range("W8").Value = 3.45
$ws.Range("X8").Value = 1.26
$ws.Range("Y8").Value = 1.39
$ws.Range("Z8").Value = 2.8
$ws.Range("AA8").Value = 1.88
$ws.Range("AB8").Value = 1.83
$ws.Range("AC8").Value = 9.75
$ws.Range("AD8").Value = 19.5
$ws.Range("AE8").Value = 13
$ws.Range("AG8").Value = 37
$ws.Range("AH8").Value = 45
$ws.Range("AI8").Value = 6.7
$ws.Range("AK8").Value = 16
$ws.Range("AL8").Value = 80
$ws.Range("AM8").Value = 700
$ws.Range("AN8").Value = 6.6
$ws.Range("AO8").Value = 8.75
$ws.Range("AR8").Value = 16
$ws.Range("AS8").Value = 30

# Row 9
$ws.Range("G9").Value = 10.25
$ws.Range("I9").Value = 1.22
$ws.Range("J9").Value = 7.6
$ws.Range("O9").Value = 1.1
$ws.Range("P9").Value = 5.9
$ws.Range("S9").Value = 1.32
$ws.Range("T9").Value = 3.1
$ws.Range("W9").Value = 1.8
$ws.Range("X9").Value = 1.91
$ws.Range("AA9").Value = 1.75
$ws.Range("AB9").Value = 1.98
$ws.Range("AD9").Value = 90
$ws.Range("AI9").Value = 27
$ws.Range("AJ9").Value = 13.5
